# Actualización automática 2025-08-14 15:45:08
$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M26").Value = 56.35
$ws1.Range("O26").Value = 1.73
$ws1.Range("P26").Value = 17.77
$ws1.Range("M55").Value = 288.72
$ws1.Range("M57").Value = "11 de 55"
$ws1.Range("O57").Value = "4 de 55"
$ws1.Range("P57").Value = "1 de 55"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F26").Value = 75.84999999999999
$ws2.Range("F55").Value = 926.14
$ws2.Range("F57").Value = 34098.1

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D10").Value = 17.77
$ws3.Range("E10").Value = 1282.73
$ws3.Range("F10").Value = 0.0136639753940792

$ws3.Range("D16").Value = 17798.14
$ws3.Range("E16").Value = 38261.56
$ws3.Range("F16").Value = 0.317485466386727

$ws3.Range("D18").Value = 1933.91
$ws3.Range("E18").Value = 1266.09
$ws3.Range("F18").Value = 0.6043468750000001

$ws3.Range("D19").Value = 34098.1
$ws3.Range("E19").Value = 83341.59064517915
$ws3.Range("F19").Value = 0.2903456217627538
